# Append a new block of 5 test-result rows (rows 67-71) to the TestResults
# sheet, replicating the existing repeating pattern of:
#   testPostVideo / testGetVideoById / testPutVideo / testDeleteVideo /
#   testInvalidPostVideo
# with a brand new timestamped "Bad Request" response captured for the
# latest testInvalidPostVideo run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestampResponse = '{"timestamp":"2025-01-02T13:08:17.564+00:00","status":400,"error":"Bad Request","path":"/api/videogame"}'

# Row 67 - testPostVideo
$ws.Range("A67").Value = "testPostVideo"
$ws.Range("B67").Value = "Passed"
$ws.Range("C67").Value = '{"category": "Platform","name": "Mario","rating": "Mature","releaseDate": "2012-05-04","reviewScore": 89,"id": "147"}'
$ws.Range("D67").Value = '{"id":0,"name":"Mario","releaseDate":"2012-05-04","reviewScore":89,"category":"Platform","rating":"Mature"}'

# Row 68 - testGetVideoById
$ws.Range("A68").Value = "testGetVideoById"
$ws.Range("B68").Value = "Passed"
$ws.Range("C68").Value = "GET https://www.videogamedb.uk:443/api/videogame/1"
$ws.Range("D68").Value = '{"id":1,"name":"Resident Evil 4","releaseDate":"2005-10-01 23:59:59","reviewScore":85,"category":"Shooter","rating":"Universal"}'

# Row 69 - testPutVideo
$ws.Range("A69").Value = "testPutVideo"
$ws.Range("B69").Value = "Passed"
$ws.Range("C69").Value = '{"category": "Platform","name": "Mario","rating": "Mature","releaseDate": "2012-05-04","reviewScore": 89,"id": "147"}'
$ws.Range("D69").Value = '{"id":1,"name":"Mario","releaseDate":"2012-05-04","reviewScore":89,"category":"Platform","rating":"Mature"}'

# Row 70 - testDeleteVideo
$ws.Range("A70").Value = "testDeleteVideo"
$ws.Range("B70").Value = "Passed"
$ws.Range("C70").Value = "DELETE https://www.videogamedb.uk:443/api/videogame/1"
$ws.Range("D70").Value = "Video game deleted"

# Row 71 - testInvalidPostVideo (new response captured with this run)
$ws.Range("A71").Value = "testInvalidPostVideo"
$ws.Range("B71").Value = "Passed"
$ws.Range("C71").Value = '{"category": "Platform","name": "InvalidGame","rating": "Everyone","releaseDate": "invalid-date","reviewScore": "invalid-score","id": "123"}'
$ws.Range("D71").Value = $newTimestampResponse
